$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D:D").Insert()

Write-Output "Inserted column"
